$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "Black Paint for Enclosure" purchase (row 39) to "Black Paint for Base"
$ws.Range("C39").Value2 = "Black Paint for Base"

# 2. Insert a new row at 40 (copy formatting/formula pattern from row 39, which shifts
#    the Subtotal/Shipping/Total rows from 41-43 down to 42-44, and the trailing blank
#    placeholder rows from 46-52 down to 47-53)
$ws.Rows("39").Copy()
$ws.Rows("40").Insert()

# 3. Populate the new row 40 with the "Frame mat" purchase details
$ws.Range("B40").Value2 = "Frame mat"
$ws.Range("C40").Value2 = "Frame mat for poster"
$ws.Range("D40").Value2 = "Michael's"
$ws.Range("E40").Value2 = 17.06
$ws.Range("F40").Value2 = 1
$ws.Range("G40").Formula = "=F40*E40"
$ws.Range("H40").Value2 = 0

# 4. Fix up the Subtotal / Shipping / Total formulas, now on rows 42-44, to include row 40
$ws.Range("G42").Formula = "=SUM(G3:G40)"
$ws.Range("G43").Formula = "=SUM(H3:H40)"
$ws.Range("G44").Formula = "=SUM(G42:H43)"

# 5. Restore the view state (scroll position / selection) to match the edited sheet
$ws.Range("G44:H44").Select()
$excel.ActiveWindow.ScrollRow = 31
